$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83; this shifts existing rows 83..131 down to 84..132
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new data entry
$ws.Cells.Item(83, 1).Value = 4
$ws.Cells.Item(83, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(83, 3).Value = "Los Lagos"
$ws.Cells.Item(83, 4).Value = 44438
$ws.Cells.Item(83, 5).Value = 10
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100108
$ws.Cells.Item(83, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(83, 9).Value = 100108005
$ws.Cells.Item(83, 10).Value = "Piña"
$ws.Cells.Item(83, 11).Value = "Caramelo"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 80
$ws.Cells.Item(83, 14).Value = 21000
$ws.Cells.Item(83, 15).Value = 21000
$ws.Cells.Item(83, 16).Value = 21000
$ws.Cells.Item(83, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(83, 18).Value = "Ecuador"
$ws.Cells.Item(83, 19).Value = 1750
$ws.Cells.Item(83, 20).Value = 12
